$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5695
$ws.Range("A2").Value = 5697.999999999956
$ws.Range("A3").Value = 5702
$ws.Range("A4").Value = 5705
$ws.Range("A5").Value = 5700
$ws.Range("A6").Value = 5691
$ws.Range("A7").Value = 5700.999999999956
$ws.Range("A8").Value = 5705.999999999985
